$d = $word.ActiveDocument

# Locate the two inline pictures we need by their known size (in points,
# derived from the wp:extent cx/cy EMU values in the target XML) rather
# than a hard-coded collection index, so the script is tied to the actual
# content being edited.
#   State Space Model picture : cx=4210638 cy=838317 EMU -> 331.55 x 66.01 pt
#   Parameter Estimation pic  : cx=5676900 cy=3114675 EMU -> 447.00 x 245.25 pt
$stateSpaceShape = $null
$paramShape = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $s = $d.InlineShapes.Item($i)
    if ([Math]::Abs($s.Width - 331.55) -lt 0.1 -and [Math]::Abs($s.Height - 66.01) -lt 0.1) {
        $stateSpaceShape = $s
    }
    if ([Math]::Abs($s.Width - 447.0) -lt 0.1 -and [Math]::Abs($s.Height - 245.25) -lt 0.1) {
        $paramShape = $s
    }
}

# --- Change 1: mark the run holding the "State Space Model" picture
#     (the inline drawing with wp14:anchorId="1826BD0D") as <w:noProof/> ---
$stateSpaceShape.Range.NoProofing = 1

# --- Change 2: the paragraph holding the "Parameter Estimation" picture
#     (inline drawing with wp14:anchorId="4B3F9206", which also carries the
#     _GoBack bookmark) needs to be followed by two new paragraphs: one
#     containing a single space, and one containing "Still more to come".
#     The bookmark must stay attached to the picture's paragraph. ---

# The paragraph immediately following the picture's paragraph (already
# exists in the document as an empty paragraph).
$followingParaIndex = $paramShape.Range.Paragraphs.Item(1).Index + 1
$followingPara = $d.Paragraphs.Item($followingParaIndex)

# Insert a new (empty) paragraph right before it -- this lands right after
# the picture's paragraph, keeping the bookmark where it belongs.
$insertPoint1 = $d.Range($followingPara.Range.Start, $followingPara.Range.Start)
$insertPoint1.InsertParagraphBefore()

# Fill that brand new paragraph with a single space.
$spacePara = $d.Paragraphs.Item($followingParaIndex)
$spaceFillPoint = $d.Range($spacePara.Range.Start, $spacePara.Range.Start)
$spaceFillPoint.InsertAfter(" ")

# Insert another new (empty) paragraph right before the original
# following paragraph (now shifted one further along).
$origFollowingPara = $d.Paragraphs.Item($followingParaIndex + 1)
$insertPoint2 = $d.Range($origFollowingPara.Range.Start, $origFollowingPara.Range.Start)
$insertPoint2.InsertParagraphBefore()

# Fill that paragraph with the new sentence.
$textPara = $d.Paragraphs.Item($followingParaIndex + 1)
$textFillPoint = $d.Range($textPara.Range.Start, $textPara.Range.Start)
$textFillPoint.InsertAfter("Still more to come")
